$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the stray empty inline-string cell at G301 (no-jail dialog fix)
$ws.Cells.Item(301, 7).ClearContents()

# Ensure text columns (H/I carry "$ N" literal strings, not currency numbers)
$hiRange = $ws.Range("H303:I318")
$hiRange.NumberFormat = "@"

# Append 16 new case rows (21CRB01459 / Hemmeter) for rows 303-318
$r = 303
$ws.Cells.Item($r, 1).Value = "21CRB01459"
$ws.Cells.Item($r, 2).Value = "Hemmeter"
$ws.Cells.Item($r, 3).Value = "POSS MARIHUANA DRUG PARAPHERNALIA"
$ws.Cells.Item($r, 4).Value = "2925.141C"
$ws.Cells.Item($r, 5).Value = "MM"
$ws.Cells.Item($r, 6).Value = "Guilty"
$ws.Cells.Item($r, 7).Value = "Guilty"
$ws.Cells.Item($r, 8).Value = "$ 100"
$ws.Cells.Item($r, 9).Value = "$ 0"

$r = 304
$ws.Cells.Item($r, 1).Value = "21CRB01459"
$ws.Cells.Item($r, 2).Value = "Hemmeter"
$ws.Cells.Item($r, 3).Value = "POSS MARIHUANA DRUG PARAPHERNALIA"
$ws.Cells.Item($r, 4).Value = "2925.141C"
$ws.Cells.Item($r, 5).Value = "MM"
$ws.Cells.Item($r, 6).Value = "Guilty"
$ws.Cells.Item($r, 7).Value = "Guilty"
$ws.Cells.Item($r, 8).Value = "$ 100"
$ws.Cells.Item($r, 9).Value = "$ 0"

$r = 305
$ws.Cells.Item($r, 1).Value = "21CRB01459"
$ws.Cells.Item($r, 2).Value = "Hemmeter"
$ws.Cells.Item($r, 3).Value = "POSS MARIHUANA DRUG PARAPHERNALIA"
$ws.Cells.Item($r, 4).Value = "2925.141C"
$ws.Cells.Item($r, 5).Value = "MM"
$ws.Cells.Item($r, 6).Value = "Guilty"
$ws.Cells.Item($r, 7).Value = "Guilty"
$ws.Cells.Item($r, 8).Value = "$ 100"
$ws.Cells.Item($r, 9).Value = "$ 0"

$r = 306
$ws.Cells.Item($r, 1).Value = "21CRB01459"
$ws.Cells.Item($r, 2).Value = "Hemmeter"
$ws.Cells.Item($r, 3).Value = "POSS MARIHUANA DRUG PARAPHERNALIA"
$ws.Cells.Item($r, 4).Value = "2925.141C"
$ws.Cells.Item($r, 5).Value = "MM"
$ws.Cells.Item($r, 6).Value = "Guilty"
$ws.Cells.Item($r, 7).Value = "Guilty"
$ws.Cells.Item($r, 8).Value = "$ 100"
$ws.Cells.Item($r, 9).Value = "$ 0"

$r = 307
$ws.Cells.Item($r, 1).Value = "21CRB01459"
$ws.Cells.Item($r, 2).Value = "Hemmeter"
$ws.Cells.Item($r, 3).Value = "POSS MARIHUANA DRUG PARAPHERNALIA"
$ws.Cells.Item($r, 4).Value = "2925.141C"
$ws.Cells.Item($r, 5).Value = "MM"
$ws.Cells.Item($r, 6).Value = "Guilty"
$ws.Cells.Item($r, 7).Value = "Guilty"
$ws.Cells.Item($r, 8).Value = "$ 100"
$ws.Cells.Item($r, 9).Value = "$ 0"

$r = 308
$ws.Cells.Item($r, 1).Value = "21CRB01459"
$ws.Cells.Item($r, 2).Value = "Hemmeter"
$ws.Cells.Item($r, 3).Value = "POSS MARIHUANA DRUG PARAPHERNALIA"
$ws.Cells.Item($r, 4).Value = "2925.141C"
$ws.Cells.Item($r, 5).Value = "MM"
$ws.Cells.Item($r, 6).Value = "Guilty"
$ws.Cells.Item($r, 7).Value = "Guilty"
$ws.Cells.Item($r, 8).Value = "$ 100"
$ws.Cells.Item($r, 9).Value = "$ 0"

$r = 309
$ws.Cells.Item($r, 1).Value = "21CRB01459"
$ws.Cells.Item($r, 2).Value = "Hemmeter"
$ws.Cells.Item($r, 3).Value = "POSS MARIHUANA DRUG PARAPHERNALIA"
$ws.Cells.Item($r, 4).Value = "2925.141C"
$ws.Cells.Item($r, 5).Value = "MM"
$ws.Cells.Item($r, 6).Value = "Guilty"
$ws.Cells.Item($r, 7).Value = "Guilty"
$ws.Cells.Item($r, 8).Value = "$ 0"
$ws.Cells.Item($r, 9).Value = "$ 0"

$r = 310
$ws.Cells.Item($r, 1).Value = "21CRB01459"
$ws.Cells.Item($r, 2).Value = "Hemmeter"
$ws.Cells.Item($r, 3).Value = "POSS MARIHUANA DRUG PARAPHERNALIA"
$ws.Cells.Item($r, 4).Value = "2925.141C"
$ws.Cells.Item($r, 5).Value = "MM"
$ws.Cells.Item($r, 6).Value = "Guilty"
$ws.Cells.Item($r, 7).Value = "Guilty"
$ws.Cells.Item($r, 8).Value = "$ 0"
$ws.Cells.Item($r, 9).Value = "$ 0"

$r = 311
$ws.Cells.Item($r, 1).Value = "21CRB01459"
$ws.Cells.Item($r, 2).Value = "Hemmeter"
$ws.Cells.Item($r, 3).Value = "POSS MARIHUANA DRUG PARAPHERNALIA"
$ws.Cells.Item($r, 4).Value = "2925.141C"
$ws.Cells.Item($r, 5).Value = "MM"
$ws.Cells.Item($r, 6).Value = "Guilty"
$ws.Cells.Item($r, 7).Value = "Guilty"
$ws.Cells.Item($r, 8).Value = "$ 0"
$ws.Cells.Item($r, 9).Value = "$ 0"

$r = 312
$ws.Cells.Item($r, 1).Value = "21CRB01459"
$ws.Cells.Item($r, 2).Value = "Hemmeter"
$ws.Cells.Item($r, 3).Value = "POSS MARIHUANA DRUG PARAPHERNALIA"
$ws.Cells.Item($r, 4).Value = "2925.141C"
$ws.Cells.Item($r, 5).Value = "MM"
$ws.Cells.Item($r, 6).Value = "Guilty"
$ws.Cells.Item($r, 7).Value = "Guilty"
$ws.Cells.Item($r, 8).Value = "$ 0"
$ws.Cells.Item($r, 9).Value = "$ 0"

$r = 313
$ws.Cells.Item($r, 1).Value = "21CRB01459"
$ws.Cells.Item($r, 2).Value = "Hemmeter"
$ws.Cells.Item($r, 3).Value = "POSS MARIHUANA DRUG PARAPHERNALIA"
$ws.Cells.Item($r, 4).Value = "2925.141C"
$ws.Cells.Item($r, 5).Value = "MM"
$ws.Cells.Item($r, 6).Value = "No Contest"
$ws.Cells.Item($r, 7).Value = "Guilty"
$ws.Cells.Item($r, 8).Value = "$ 0"
$ws.Cells.Item($r, 9).Value = "$ 0"

$r = 314
$ws.Cells.Item($r, 1).Value = "21CRB01459"
$ws.Cells.Item($r, 2).Value = "Hemmeter"
$ws.Cells.Item($r, 3).Value = "POSS MARIHUANA DRUG PARAPHERNALIA"
$ws.Cells.Item($r, 4).Value = "2925.141C"
$ws.Cells.Item($r, 5).Value = "MM"
$ws.Cells.Item($r, 6).Value = "No Contest"
$ws.Cells.Item($r, 7).Value = "Guilty"
$ws.Cells.Item($r, 8).Value = "$ 0"
$ws.Cells.Item($r, 9).Value = "$ 0"

$r = 315
$ws.Cells.Item($r, 1).Value = "21CRB01459"
$ws.Cells.Item($r, 2).Value = "Hemmeter"
$ws.Cells.Item($r, 3).Value = "POSS MARIHUANA DRUG PARAPHERNALIA"
$ws.Cells.Item($r, 4).Value = "2925.141C"
$ws.Cells.Item($r, 5).Value = "MM"
$ws.Cells.Item($r, 6).Value = "No Contest"
$ws.Cells.Item($r, 7).Value = "Guilty"
$ws.Cells.Item($r, 8).Value = "$ 0"
$ws.Cells.Item($r, 9).Value = "$ 0"

$r = 316
$ws.Cells.Item($r, 1).Value = "21CRB01459"
$ws.Cells.Item($r, 2).Value = "Hemmeter"
$ws.Cells.Item($r, 3).Value = "POSS MARIHUANA DRUG PARAPHERNALIA"
$ws.Cells.Item($r, 4).Value = "2925.141C"
$ws.Cells.Item($r, 5).Value = "MM"
$ws.Cells.Item($r, 6).Value = "Guilty"
$ws.Cells.Item($r, 7).Value = "Guilty"
$ws.Cells.Item($r, 8).Value = "$ 100"
$ws.Cells.Item($r, 9).Value = "$ 0"

$r = 317
$ws.Cells.Item($r, 1).Value = "21CRB01459"
$ws.Cells.Item($r, 2).Value = "Hemmeter"
$ws.Cells.Item($r, 3).Value = "POSS MARIHUANA DRUG PARAPHERNALIA"
$ws.Cells.Item($r, 4).Value = "2925.141C"
$ws.Cells.Item($r, 5).Value = "MM"
$ws.Cells.Item($r, 6).Value = "Guilty"
$ws.Cells.Item($r, 7).Value = "Guilty"
$ws.Cells.Item($r, 8).Value = "$ 100"
$ws.Cells.Item($r, 9).Value = "$ 0"

$r = 318
$ws.Cells.Item($r, 1).Value = "21CRB01459"
$ws.Cells.Item($r, 2).Value = "Hemmeter"
$ws.Cells.Item($r, 3).Value = "POSS MARIHUANA DRUG PARAPHERNALIA"
$ws.Cells.Item($r, 4).Value = "2925.141C"
$ws.Cells.Item($r, 5).Value = "MM"
$ws.Cells.Item($r, 6).Value = "Guilty"
$ws.Cells.Item($r, 7).Value = "Guilty"
$ws.Cells.Item($r, 8).Value = "$ 0"
$ws.Cells.Item($r, 9).Value = "$ 0"
